$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the style of the existing header row (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Fill the new "Save" column data
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
